$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Refresh cryptocurrency Price/Volume(1h) columns with latest scraped values
# (GitHub Actions scheduled update, 2024-05-25). A couple of coins also swapped
# rank position (rows 31-34), so Coin/Link/Price/Volume are rewritten there too.

$ws.Cells.Item(2, 4).Value = '69.106.98'
$ws.Cells.Item(2, 5).Value = '  +0.23%  '

$ws.Cells.Item(3, 4).Value = '3.748.46'
$ws.Cells.Item(3, 5).Value = '  +0.20%  '

$ws.Cells.Item(4, 5).Value = '  +0.00%  '

$ws.Cells.Item(5, 4).Value = "'601.42"
$ws.Cells.Item(5, 5).Value = '  -0.01%  '

$ws.Cells.Item(6, 4).Value = "'166.95"
$ws.Cells.Item(6, 5).Value = '  -0.61%  '

$ws.Cells.Item(7, 4).Value = '3.746.01'
$ws.Cells.Item(7, 5).Value = '  +0.17%  '

$ws.Cells.Item(8, 5).Value = '  +0.00%  '

$ws.Cells.Item(9, 5).Value = '  +1.16%  '

$ws.Cells.Item(10, 5).Value = '  +2.58%  '

$ws.Cells.Item(11, 4).Value = "'6.38"
$ws.Cells.Item(11, 5).Value = '  +0.65%  '

$ws.Cells.Item(12, 5).Value = '  +0.06%  '

$ws.Cells.Item(13, 4).Value = "'37.94"
$ws.Cells.Item(13, 5).Value = '  -0.66%  '

$ws.Cells.Item(14, 5).Value = '  +1.44%  '

$ws.Cells.Item(15, 4).Value = '4.374.26'

$ws.Cells.Item(16, 4).Value = '3.749.08'
$ws.Cells.Item(16, 5).Value = '  +0.27%  '

$ws.Cells.Item(17, 4).Value = '69.101.34'
$ws.Cells.Item(17, 5).Value = '  +0.27%  '

$ws.Cells.Item(18, 4).Value = "'7.35"
$ws.Cells.Item(18, 5).Value = '  +1.50%  '

$ws.Cells.Item(19, 4).Value = "'17.40"
$ws.Cells.Item(19, 5).Value = '  +1.02%  '

$ws.Cells.Item(20, 5).Value = '  -1.69%  '

$ws.Cells.Item(21, 5).Value = '  +8.48%  '

$ws.Cells.Item(22, 4).Value = "'492.24"
$ws.Cells.Item(22, 5).Value = '  -0.87%  '

$ws.Cells.Item(23, 4).Value = "'0.727"
$ws.Cells.Item(23, 5).Value = '  +0.62%  '

$ws.Cells.Item(24, 5).Value = '  +7.38%  '

$ws.Cells.Item(25, 4).Value = "'84.85"
$ws.Cells.Item(25, 5).Value = '  -0.50%  '

$ws.Cells.Item(26, 4).Value = "'2.29"
$ws.Cells.Item(26, 5).Value = '  -0.30%  '

$ws.Cells.Item(27, 4).Value = "'12.26"
$ws.Cells.Item(27, 5).Value = '  -0.13%  '

$ws.Cells.Item(28, 4).Value = "'10.05"
$ws.Cells.Item(28, 5).Value = '  -0.75%  '

$ws.Cells.Item(29, 5).Value = '  -0.09%  '

$ws.Cells.Item(30, 5).Value = '  +0.94%  '

$ws.Cells.Item(31, 2).Value = 'ImmutableX'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(31, 4).Value = "'2.47"
$ws.Cells.Item(31, 5).Value = '  +2.12%  '

$ws.Cells.Item(32, 2).Value = 'NEARProtocol'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(32, 4).Value = "'8.11"
$ws.Cells.Item(32, 5).Value = '  +1.97%  '

$ws.Cells.Item(33, 2).Value = 'WrappedeETH'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Cells.Item(33, 4).Value = '3.895.19'
$ws.Cells.Item(33, 5).Value = '  +0.07%  '

$ws.Cells.Item(34, 2).Value = 'EthereumClassic'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(34, 4).Value = "'31.45"
$ws.Cells.Item(34, 5).Value = '  -0.92%  '

$ws.Cells.Item(35, 4).Value = '3.681.40'
$ws.Cells.Item(35, 5).Value = '  +0.25%  '

$ws.Cells.Item(36, 5).Value = '  -0.22%  '

$ws.Cells.Item(37, 5).Value = '  +0.01%  '

$ws.Cells.Item(39, 5).Value = '  -0.01%  '

$ws.Cells.Item(40, 4).Value = "'0.138"
$ws.Cells.Item(40, 5).Value = '  +2.91%  '

$ws.Cells.Item(41, 5).Value = '  +0.35%  '

$ws.Cells.Item(42, 5).Value = '  +5.38%  '

$ws.Cells.Item(43, 5).Value = '  -0.75%  '

$ws.Cells.Item(44, 4).Value = "'425.55"
$ws.Cells.Item(44, 5).Value = '  -2.47%  '

$ws.Cells.Item(45, 5).Value = '  -0.56%  '

$ws.Cells.Item(46, 4).Value = "'8.46"
$ws.Cells.Item(46, 5).Value = '  +0.61%  '

$ws.Cells.Item(48, 5).Value = '  -0.71%  '

$ws.Cells.Item(49, 4).Value = "'141.28"
$ws.Cells.Item(49, 5).Value = '  -0.56%  '

$ws.Cells.Item(50, 4).Value = '2.789.37'
$ws.Cells.Item(50, 5).Value = '  +1.50%  '

$ws.Cells.Item(51, 5).Value = '  +0.38%  '
